$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("dSF") updated values per row
$values = @{
    2  = 2
    3  = -2
    4  = -1
    5  = -3
    6  = 1
    7  = 0
    8  = -5
    9  = -1
    10 = 1
    11 = 1
    12 = 1
    13 = -1
    14 = -12
    16 = -4
    17 = 3
    18 = -1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
